$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata": bump Version, Date and Contact text (new term release,
# now also carrying a proper ContactDetail display instead of the old
# "No display for ContactDetail" placeholder).
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# ---------------------------------------------------------------------------
# Sheet "Include from FSIII": a new Concept row is inserted before the
# trailing "System URI" row. The old row 13 (previously an empty spare
# Concept/Description pair) now gets the new concept code, the
# "System URI" / "urn:oid:..." pair moves down to row 15, and row 14
# becomes the new empty spacer row.
# ---------------------------------------------------------------------------
$inc = $wb.Worksheets.Item("Include from FSIII")

# Make room: move the "System URI" row down to row 15, carrying the regular
# data-row formatting (style) along with it (format-only paste from an
# existing plain data row keeps the same shared cell style index).
$inc.Range("A12:B12").Copy()
$inc.Range("A15:B15").PasteSpecial(-4122)
$inc.Range("A15").Value = "System URI"
$inc.Range("B15").Value = "urn:oid:1.2.208.176.2.21"

# Row 14 becomes the new empty spacer row (same style, no content).
$inc.Range("A12:B12").Copy()
$inc.Range("A14:B14").PasteSpecial(-4122)
$inc.Range("A14").ClearContents()
$inc.Range("B14").ClearContents()

# Row 13 gets the new Concept code; its Description cell stays empty.
$inc.Range("A13").Value = "aec684bd-c2ea-4ff0-8eb7-6d2cf67fb863"
$inc.Range("B13").ClearContents()
